# Apply updated crypto price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.296.27'
$ws.Range("E2").Value = '  -5.68%  '
$ws.Range("D3").Value = '1.670.57'
$ws.Range("E3").Value = '  -3.83%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "`'217.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.78%  '
$ws.Range("D6").Value = "`'0.5085"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -11.78%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = "`'0.2657"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").Value = "`'0.06352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.04%  '
$ws.Range("D10").Value = "`'21.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.52%  '
$ws.Range("D11").Value = "`'0.07369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D12").Value = '1.668.20'
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").Value = "`'4.551"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("D14").Value = "`'0.5817"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.27%  '
$ws.Range("D15").Value = '1.896.34'
$ws.Range("E15").Value = '  -3.95%  '
$ws.Range("D16").Value = "`'0.000008522"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("D17").Value = "`'64.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -12.94%  '
$ws.Range("D18").Value = '26.357.60'
$ws.Range("E18").Value = '  -5.40%  '
$ws.Range("D19").Value = "`'4.939"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.95%  '
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = "`'10.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.71%  '
$ws.Range("D22").Value = "`'188.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.75%  '
$ws.Range("D23").Value = "`'6.201"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.28%  '
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Value = "`'143.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.46%  '
$ws.Range("D26").Value = "`'7.675"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.14%  '
$ws.Range("D27").Value = "`'0.1178"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.27%  '
$ws.Range("D28").Value = "`'15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.21%  '
$ws.Range("D29").Value = "`'0.05878"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.83%  '
$ws.Range("D30").Value = "`'1.268"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.47%  '
$ws.Range("D31").Value = "`'1.321"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.04%  '
$ws.Range("D32").Value = "`'3.532"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.19%  '
$ws.Range("D33").Value = "`'3.514"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.03%  '
$ws.Range("D34").Value = "`'1.641"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("D35").Value = "`'1.013"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.02%  '
$ws.Range("D36").Value = "`'0.6004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.83%  '
$ws.Range("D37").Value = "`'2.355"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.40%  '
$ws.Range("D38").Value = "`'2.647"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("D39").Value = "`'0.01615"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("D40").Value = "`'6.015"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").Value = '1.074.03'
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("D42").Value = "`'0.8662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").Value = "`'99.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").Value = '1.819.29'
$ws.Range("E45").Value = '  -3.59%  '
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").Value = "`'55.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.82%  '
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = "`'8.085"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("D50").Value = "`'0.4291"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("D51").Value = "`'0.05180"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.61%  '
